$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 248
$ws1.Range("F5").Value = 264
$ws1.Range("F7").Value = 1394
$ws1.Range("F14").Value = 398
$ws1.Range("F18").Value = 259
$ws1.Range("F23").Value = 5464
$ws1.Range("F28").Value = 13870
$ws1.Range("F31").Value = 82
$ws1.Range("F33").Value = 391
$ws1.Range("F34").Value = 542
$ws1.Range("F36").Value = 82
$ws1.Range("F38").Value = 110

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 247
$ws4.Range("F5").Value = 264
$ws4.Range("F13").Value = 107
$ws4.Range("F14").Value = 398
$ws4.Range("F18").Value = 259
$ws4.Range("F26").Value = 5464
$ws4.Range("F31").Value = 13870
$ws4.Range("F34").Value = 82
$ws4.Range("F36").Value = 391
$ws4.Range("F37").Value = 542
$ws4.Range("F38").Value = 4135
$ws4.Range("F39").Value = 82
$ws4.Range("F41").Value = 110
